$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 163.9108346666667
$ws.Range("H2").Value = 491.732504
$ws.Range("I2").Value = 0.8426759240348239
$ws.Range("J2").Value = 0.8426759240348242
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.1743356666666667
$ws.Range("N2").Value = 0.523007
$ws.Range("O2").Value = 0.02303006925000699
$ws.Range("P2").Value = 0.02303006925000699
$ws.Range("Q2").Value = 28.57550463550311
$ws.Range("R2").Value = 257.179541719528
$ws.Range("S2").Value = 0.01940688488583563
$ws.Range("T2").Value = 0.01940688488583563

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 163.9108346666667
$ws.Range("H3").Value = 491.732504
$ws.Range("I3").Value = 0.8426759240348239
$ws.Range("J3").Value = 0.8426759240348242
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.401382333333334
$ws.Range("N3").Value = 7.204147000000001
$ws.Range("O3").Value = 0.3172271198994089
$ws.Range("P3").Value = 0.3172271198994089
$ws.Range("Q3").Value = 393.6125826104543
$ws.Range("R3").Value = 3542.513243494088
$ws.Range("S3").Value = 0.2673196563901403
$ws.Range("T3").Value = 0.2673196563901403

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 163.9108346666667
$ws.Range("H4").Value = 491.732504
$ws.Range("I4").Value = 0.8426759240348239
$ws.Range("J4").Value = 0.8426759240348242
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.994197
$ws.Range("N4").Value = 14.982591
$ws.Range("O4").Value = 0.6597428108505842
$ws.Range("P4").Value = 0.6597428108505842
$ws.Range("Q4").Value = 818.6029987597626
$ws.Range("R4").Value = 7367.426988837864
$ws.Range("S4").Value = 0.5559493827588481
$ws.Range("T4").Value = 0.5559493827588482

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 29.35342966666667
$ws.Range("H5").Value = 88.060289
$ws.Range("I5").Value = 0.1509078305790594
$ws.Range("J5").Value = 0.1509078305790594
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.1743356666666667
$ws.Range("N5").Value = 0.523007
$ws.Range("O5").Value = 0.02303006925000699
$ws.Range("P5").Value = 0.02303006925000699
$ws.Range("Q5").Value = 5.117349729891444
$ws.Range("R5").Value = 46.056147569023
$ws.Range("S5").Value = 0.003475417788604061
$ws.Range("T5").Value = 0.003475417788604061

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 29.35342966666667
$ws.Range("H6").Value = 88.060289
$ws.Range("I6").Value = 0.1509078305790594
$ws.Range("J6").Value = 0.1509078305790594
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.401382333333334
$ws.Range("N6").Value = 7.204147000000001
$ws.Range("O6").Value = 0.3172271198994089
$ws.Range("P6").Value = 0.3172271198994089
$ws.Range("Q6").Value = 70.4888074242759
$ws.Range("R6").Value = 634.3992668184831
$ws.Range("S6").Value = 0.04787205646486296
$ws.Range("T6").Value = 0.04787205646486296

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 29.35342966666667
$ws.Range("H7").Value = 88.060289
$ws.Range("I7").Value = 0.1509078305790594
$ws.Range("J7").Value = 0.1509078305790594
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.994197
$ws.Range("N7").Value = 14.982591
$ws.Range("O7").Value = 0.6597428108505842
$ws.Range("P7").Value = 0.6597428108505842
$ws.Range("Q7").Value = 146.5968103809777
$ws.Range("R7").Value = 1319.371293428799
$ws.Range("S7").Value = 0.09956035632559239
$ws.Range("T7").Value = 0.09956035632559239

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.248038666666667
$ws.Range("H8").Value = 3.744116
$ws.Range("I8").Value = 0.006416245386116614
$ws.Range("J8").Value = 0.006416245386116614
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.1743356666666667
$ws.Range("N8").Value = 0.523007
$ws.Range("O8").Value = 0.02303006925000699
$ws.Range("P8").Value = 0.02303006925000699
$ws.Range("Q8").Value = 0.2175776529791111
$ws.Range("R8").Value = 1.958198876812
$ws.Range("S8").Value = 0.0001477665755673034
$ws.Range("T8").Value = 0.0001477665755673035

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.248038666666667
$ws.Range("H9").Value = 3.744116
$ws.Range("I9").Value = 0.006416245386116614
$ws.Range("J9").Value = 0.006416245386116614
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.401382333333334
$ws.Range("N9").Value = 7.204147000000001
$ws.Range("O9").Value = 0.3172271198994089
$ws.Range("P9").Value = 0.3172271198994089
$ws.Range("Q9").Value = 2.997018005450223
$ws.Range("R9").Value = 26.973162049052
$ws.Range("S9").Value = 0.002035407044405644
$ws.Range("T9").Value = 0.002035407044405644

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.248038666666667
$ws.Range("H10").Value = 3.744116
$ws.Range("I10").Value = 0.006416245386116614
$ws.Range("J10").Value = 0.006416245386116614
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 4.994197
$ws.Range("N10").Value = 14.982591
$ws.Range("O10").Value = 0.6597428108505842
$ws.Range("P10").Value = 0.6597428108505842
$ws.Range("Q10").Value = 6.232950964950667
$ws.Range("R10").Value = 56.09655868455599
$ws.Range("S10").Value = 0.004233071766143666
$ws.Range("T10").Value = 0.004233071766143667
